$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "62.034.86"
$ws.Cells.Item(2, 5).Value = "  +2.44%  "
$ws.Cells.Item(3, 4).Value = "2.415.40"
$ws.Cells.Item(3, 5).Value = "  +0.10%  "
$ws.Cells.Item(4, 5).Value = "  +0.61%  "
$ws.Cells.Item(5, 4).Value = "'570.56"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.11%  "
$ws.Cells.Item(6, 4).Value = "'143.94"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +4.73%  "
$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.45%  "
$ws.Cells.Item(8, 4).Value = "'0.536"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.55%  "
$ws.Cells.Item(9, 4).Value = "2.430.26"
$ws.Cells.Item(9, 5).Value = "  +1.38%  "
$ws.Cells.Item(10, 5).Value = "  +4.42%  "
$ws.Cells.Item(11, 5).Value = "  +0.76%  "
$ws.Cells.Item(12, 4).Value = "'5.23"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +4.03%  "
$ws.Cells.Item(13, 4).Value = "'0.348"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +3.72%  "
$ws.Cells.Item(14, 4).Value = "'26.61"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +3.65%  "
$ws.Cells.Item(15, 4).Value = "'0.0000176"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +6.20%  "
$ws.Cells.Item(16, 4).Value = "2.842.52"
$ws.Cells.Item(16, 5).Value = "  +1.12%  "
$ws.Cells.Item(17, 4).Value = "61.966.11"
$ws.Cells.Item(17, 5).Value = "  +2.09%  "
$ws.Cells.Item(18, 4).Value = "2.424.65"
$ws.Cells.Item(18, 5).Value = "  +2.03%  "
$ws.Cells.Item(19, 5).Value = "  -2.62%  "
$ws.Cells.Item(20, 4).Value = "'10.80"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.48%  "
$ws.Cells.Item(21, 4).Value = "'325.64"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.16%  "
$ws.Cells.Item(22, 4).Value = "'4.13"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +2.39%  "
$ws.Cells.Item(23, 2).Value = "SuiNetwork"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(23, 4).Value = "'2.02"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +12.89%  "
$ws.Cells.Item(24, 2).Value = "Dai"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(24, 4).Value = "'1.00"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.16%  "
$ws.Cells.Item(25, 2).Value = "Litecoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(25, 4).Value = "'65.20"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.64%  "
$ws.Cells.Item(26, 2).Value = "Bittensor"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(26, 4).Value = "'611.38"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +10.89%  "
$ws.Cells.Item(27, 2).Value = "Aptos"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(27, 4).Value = "'8.40"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +4.42%  "
$ws.Cells.Item(28, 2).Value = "PEPE"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(28, 4).Value = "0.0₃0983"
$ws.Cells.Item(28, 5).Value = "  +8.42%  "
$ws.Cells.Item(29, 2).Value = "WrappedeETH"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(29, 4).Value = "2.537.86"
$ws.Cells.Item(29, 5).Value = "  +0.44%  "
$ws.Cells.Item(30, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(30, 4).Value = "'8.04"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.22%  "
$ws.Cells.Item(31, 2).Value = "Fetch.AI"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(31, 4).Value = "'1.41"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +8.97%  "
$ws.Cells.Item(32, 2).Value = "PancakeSwap"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(32, 4).Value = "'1.82"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.63%  "
$ws.Cells.Item(33, 2).Value = "Kaspa"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(33, 4).Value = "'0.135"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +3.42%  "
$ws.Cells.Item(34, 2).Value = "ImmutableX"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(34, 4).Value = "'1.47"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +4.52%  "
$ws.Cells.Item(35, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(35, 4).Value = "'0.997"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.53%  "
$ws.Cells.Item(36, 2).Value = "NEARProtocol"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(36, 4).Value = "'4.75"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +5.69%  "
$ws.Cells.Item(37, 4).Value = "'152.92"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.25%  "
$ws.Cells.Item(38, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(38, 4).Value = "'0.371"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.16%  "
$ws.Cells.Item(39, 2).Value = "RenderToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(39, 4).Value = "'5.36"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +5.89%  "
$ws.Cells.Item(40, 2).Value = "EthereumClassic"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(40, 4).Value = "'18.50"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.76%  "
$ws.Cells.Item(41, 2).Value = "dogwifhat"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(41, 4).Value = "'2.65"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +16.55%  "
$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(42, 4).Value = "'1.72"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +5.07%  "
$ws.Cells.Item(43, 2).Value = "USDe"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(43, 4).Value = "'0.999"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.05%  "
$ws.Cells.Item(44, 4).Value = "'42.15"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.94%  "
$ws.Cells.Item(45, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(45, 4).Value = "0.0₆0283"
$ws.Cells.Item(45, 5).Value = "  -1.94%  "
$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).Value = "'142.57"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.38%  "
$ws.Cells.Item(47, 2).Value = "Filecoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(47, 4).Value = "'3.57"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.36%  "
$ws.Cells.Item(48, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(48, 4).Value = "'20.20"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +6.39%  "
$ws.Cells.Item(49, 2).Value = "Mantle"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(49, 4).Value = "'0.598"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.71%  "
$ws.Cells.Item(50, 2).Value = "Hedera"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(50, 4).Value = "'0.0512"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +3.36%  "
$ws.Cells.Item(51, 2).Value = "Stellar"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(51, 4).Value = "'0.0916"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +2.55%  "
